$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.170247316360474
$ws.Range("B1").Value = 1.665533781051636
$ws.Range("C1").Value = 4.56535816192627
$ws.Range("D1").Value = 0.6432885527610779
$ws.Range("E1").Value = 0.7118012309074402
